$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.00353462736640307
$ws.Cells.Item(2, 3).Value = 0.985202492211838
$ws.Cells.Item(2, 4).Value = 0.020968128444764
$ws.Cells.Item(2, 5).Value = 0.00365444524323029
$ws.Cells.Item(2, 6).Value = 0.000299544692068057
$ws.Cells.Item(2, 7).Value = 0.990534387730649
$ws.Cells.Item(2, 8).Value = 0.000299544692068057
$ws.Cells.Item(2, 9).Value = 0.000299544692068057
$ws.Cells.Item(2, 10).Value = 0.0439731607955907
$ws.Cells.Item(2, 11).Value = 0.0000599089384136113
$ws.Cells.Item(2, 12).Value = 0.00353462736640307
$ws.Cells.Item(2, 13).Value = 0.00143781452192667
$ws.Cells.Item(2, 14).Value = 0.69260723699976
$ws.Cells.Item(2, 15).Value = 0.00101845195303139
$ws.Cells.Item(2, 16).Value = 0.00569134914929307
$ws.Cells.Item(2, 17).Value = 0.000179726815240834
$ws.Cells.Item(2, 18).Value = 0
$ws.Cells.Item(2, 19).Value = 0.986220944164869
$ws.Cells.Item(2, 20).Value = 0.0000599089384136113
$ws.Cells.Item(2, 21).Value = 0
$ws.Cells.Item(2, 22).Value = 0.972082434699257
$ws.Cells.Item(2, 23).Value = 0.0101845195303139
$ws.Cells.Item(2, 24).Value = 0.00658998322549724

$ws.Cells.Item(3, 2).Value = 0.995866283249461
$ws.Cells.Item(3, 3).Value = 0.0101246105919003
$ws.Cells.Item(3, 4).Value = 0.00832734243949197
$ws.Cells.Item(3, 5).Value = 0.000299544692068057
$ws.Cells.Item(3, 6).Value = 0.995446920680566
$ws.Cells.Item(3, 7).Value = 0.000299544692068057
$ws.Cells.Item(3, 8).Value = 0.000239635753654445
$ws.Cells.Item(3, 9).Value = 0.000119817876827223
$ws.Cells.Item(3, 10).Value = 0.000958543014617781
$ws.Cells.Item(3, 11).Value = 0.996105919003115
$ws.Cells.Item(3, 12).Value = 0.99568655643422
$ws.Cells.Item(3, 13).Value = 0.0224059429666906
$ws.Cells.Item(3, 14).Value = 0.001078360891445
$ws.Cells.Item(3, 15).Value = 0.0328900071890726
$ws.Cells.Item(3, 16).Value = 0.0307931943445962
$ws.Cells.Item(3, 17).Value = 0.00143781452192667
$ws.Cells.Item(3, 18).Value = 0.999640546369518
$ws.Cells.Item(3, 19).Value = 0.00245626647495806
$ws.Cells.Item(3, 20).Value = 0.996645099448838
$ws.Cells.Item(3, 21).Value = 0.999760364246346
$ws.Cells.Item(3, 22).Value = 0.0228852144739995
$ws.Cells.Item(3, 23).Value = 0.000119817876827223
$ws.Cells.Item(3, 24).Value = 0.000539180445722502

$ws.Cells.Item(4, 2).Value = 0.000479271507308891
$ws.Cells.Item(4, 3).Value = 0.00401389887371196
$ws.Cells.Item(4, 4).Value = 0.166127486220944
$ws.Cells.Item(4, 5).Value = 0.995746465372634
$ws.Cells.Item(4, 6).Value = 0.000179726815240834
$ws.Cells.Item(4, 7).Value = 0.00814761562425114
$ws.Cells.Item(4, 8).Value = 0.999281092739037
$ws.Cells.Item(4, 9).Value = 0.999580637431105
$ws.Cells.Item(4, 10).Value = 0.953989935298347
$ws.Cells.Item(4, 11).Value = 0.000239635753654445
$ws.Cells.Item(4, 12).Value = 0.000119817876827223
$ws.Cells.Item(4, 13).Value = 0.00653007428708363
$ws.Cells.Item(4, 14).Value = 0.300503235082674
$ws.Cells.Item(4, 15).Value = 0.0000599089384136113
$ws.Cells.Item(4, 16).Value = 0.00820752456266475
$ws.Cells.Item(4, 17).Value = 0.998142822909178
$ws.Cells.Item(4, 18).Value = 0.0000599089384136113
$ws.Cells.Item(4, 19).Value = 0.00976515696141864
$ws.Cells.Item(4, 20).Value = 0.000119817876827223
$ws.Cells.Item(4, 21).Value = 0.000119817876827223
$ws.Cells.Item(4, 22).Value = 0.00383417205847112
$ws.Cells.Item(4, 23).Value = 0.989456026839204
$ws.Cells.Item(4, 24).Value = 0.992751018451953

$ws.Cells.Item(5, 2).Value = 0.0000599089384136113
$ws.Cells.Item(5, 3).Value = 0.000359453630481668
$ws.Cells.Item(5, 4).Value = 0.799424874191229
$ws.Cells.Item(5, 5).Value = 0.000179726815240834
$ws.Cells.Item(5, 6).Value = 0.00377426312005751
$ws.Cells.Item(5, 7).Value = 0.000479271507308891
$ws.Cells.Item(5, 8).Value = 0.000179726815240834
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0.00353462736640307
$ws.Cells.Item(5, 12).Value = 0.000419362568895279
$ws.Cells.Item(5, 13).Value = 0.968607716271268
$ws.Cells.Item(5, 14).Value = 0.0000599089384136113
$ws.Cells.Item(5, 15).Value = 0.964653726335969
$ws.Cells.Item(5, 16).Value = 0.953810208483106
$ws.Cells.Item(5, 17).Value = 0.000119817876827223
$ws.Cells.Item(5, 18).Value = 0.000299544692068057
$ws.Cells.Item(5, 19).Value = 0.000539180445722502
$ws.Cells.Item(5, 20).Value = 0.00305535585909418
$ws.Cells.Item(5, 21).Value = 0.000119817876827223
$ws.Cells.Item(5, 22).Value = 0.000838725137790558
$ws.Cells.Item(5, 23).Value = 0.0000599089384136113
$ws.Cells.Item(5, 24).Value = 0.000119817876827223

$wb.Save()
